$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") ---
# Copy the formatting (bold font, thin borders, center/top alignment) that the
# existing header cells (e.g. H1) already use, then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- New data columns I and J for rows 2-52 ---
$iValues = @(8,7,8,8,7,7,7,8,8,8,8,8,9,8,8,8,9,8,8,9,10,9,9,8,7,8,6,11,8,8,8,9,8,9,7,9,9,7,7,7,7,9,5,7,4,4,3,7,7,7,5)
$jValues = @(8,8,8,8,8,7,7,8,8,8,8,8,9,8,8,9,9,8,8,9,10,9,10,8,8,8,6,11,8,8,8,9,8,9,7,9,9,7,8,7,7,9,5,7,4,4,3,7,7,7,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
